$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (CO 3013696596) becomes the new fully-populated data row, relocated to row 8.
$ws.Range("A11:O11").Cut($ws.Range("A8:O8"))

# New single-cell row: TS2 marker on row 7.
$ws.Range("A7").Value = "TS2"

# New partially-populated row 9 (CO 3013696923, product TA5ACMFTWLC only).
$ws.Range("A9").Value = "3013696923"
$ws.Range("B9").Value = "TA5ACMFTWLC"

# Row 2 keeps only its CO-number cell; old detail cells B2:O2 are cleared out.
$ws.Range("B2:O2").ClearContents()
$ws.Range("A2").Value = "3013996644"

# Drop the now-unused old rows (old row10 data, and the row11 remnants left by Cut).
$ws.Rows("10:11").Delete()

$ws.Range("A2").Select()
